$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns at R:T (shifting existing R:AE -> U:AH), preserving styles via xlShiftToRight.
$ws.Range("R1:T2").Insert(-4161)

# New header labels for the inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data values for row 2 in the inserted columns.
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Normalize casing of existing text values in row 2 (D2:J2), and fix H2's content.
$ws.Range("D2").Value = "considered"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "not considered"
$ws.Range("G2").Value = "very important"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"
